$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.191716432571411
$ws.Range("B1").Value = 3.010075807571411
$ws.Range("C1").Value = 5.442152976989746
$ws.Range("D1").Value = 2.291221380233765
$ws.Range("E1").Value = 1.398504853248596
